# Auto-generated Excel COM-interop script to update crypto price/volume data
# (Tue Feb 13 20:23:45 UTC 2024 GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '49.203.94'
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -1.24%  '
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.619.73'
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -0.01%  '
$cell.ClearFormats()

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.07%  '
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '112.10'
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  +1.89%  '
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '323.03'
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -1.31%  '
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  -1.32%  '
$cell.ClearFormats()

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.542'
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -3.19%  '
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '39.77'
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -1.58%  '
$cell.ClearFormats()

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '19.72'
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -4.68%  '
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.0810'
$cell.ClearFormats()

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.127'
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  +1.09%  '
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '7.26'
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.029.43'
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -0.02%  '
$cell.ClearFormats()

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.630.14'
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  +0.57%  '
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.856'
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  -1.81%  '
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '49.149.42'
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  -1.22%  '
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  -1.52%  '
$cell.ClearFormats()

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '12.90'
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -3.54%  '
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -2.21%  '
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0943'
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  -1.14%  '
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '269.20'
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  -3.34%  '
$cell.ClearFormats()

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '68.53'
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -5.52%  '
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  -2.06%  '
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '26.13'
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  -1.42%  '
$cell.ClearFormats()

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  +0.02%  '
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '10.31'
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  +3.40%  '
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -0.34%  '
$cell.ClearFormats()

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.139'
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -3.97%  '
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -4.97%  '
$cell.ClearFormats()

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '49.57'
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  -0.40%  '
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '5.46'
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  +0.37%  '
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '0.0809'
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  +2.44%  '
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -0.32%  '
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '18.99'
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -3.86%  '
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  +3.75%  '
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.04'
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -0.71%  '
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '3.13'
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  +1.39%  '
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '126.73'
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  +2.46%  '
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -1.66%  '
$cell.ClearFormats()

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '22.19'
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -2.12%  '
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0321'
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  +1.86%  '
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -4.27%  '
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '2.059.14'
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  +0.46%  '
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  +7.36%  '
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '3.21'
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -4.17%  '
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -5.57%  '
$cell.ClearFormats()

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '8.91'
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -1.48%  '
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '58.92'
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  +2.01%  '
$cell.ClearFormats()

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '5.19'
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -3.03%  '
$cell.ClearFormats()

